$d = $word.ActiveDocument

# Replace version "4.2.0" -> "4.2.1" in all occurrences (textbox + VML fallback)
$d.Content.Find.Execute(".0", $true, $false, $false, $false, $false, $true, 1, $false, ".1", 2) | Out-Null

# Replace "January 2024" -> "March 2024" in all occurrences
$d.Content.Find.Execute("January 2024", $true, $false, $false, $false, $false, $true, 1, $false, "March 2024", 2) | Out-Null
